$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A header/labels (folder / company naming instead of person name)
$ws.Range("A1").Value = "Yrityksen nimi"
$ws.Range("A2").Value = "RoboCamp demoyritys 1"
$ws.Range("A3").Value = "Ohjelmistorobotti"

# Select A3 to match resulting cursor position
$ws.Range("A3").Select()
